$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 646, shifting existing rows 646:697 down to 647:698
$ws.Rows.Item(646).Insert()

# Populate the newly inserted row 646 with the new record's data.
# Columns A, B, C, E, F, G, N, O, Q, R mirror the surrounding rows (unchanged).
$ws.Range("A646").Value = 4
$ws.Range("B646").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C646").Value = "Los Lagos"
$ws.Range("D646").Value = 45021
$ws.Range("E646").Value = 10
$ws.Range("F646").Value = 100112006
$ws.Range("G646").Value = "Repollo"
$ws.Range("H646").Value = "Crespo record"
$ws.Range("I646").Value = "Primera"
$ws.Range("J646").Value = 100
$ws.Range("K646").Value = 2000
$ws.Range("L646").Value = 2000
$ws.Range("M646").Value = 2000
$ws.Range("N646").Value = "$/unidad"
$ws.Range("O646").Value = "Región Metropolitana"
$ws.Range("P646").Value = 2000
$ws.Range("Q646").Value = 1
$ws.Range("R646").Value = "Hortaliza"
